$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 44279
$ws.Range("J2").Value = 60
$ws.Range("K2").Value = 7500
$ws.Range("L2").Value = 8000
$ws.Range("M2").Value = 7750
$ws.Range("P2").Value = 775
$ws.Range("D3").Value = 44672
$ws.Range("J3").Value = 220
$ws.Range("M3").Value = 10455
$ws.Range("P3").Value = 1046
$ws.Range("D4").Value = 44320
$ws.Range("J4").Value = 50
$ws.Range("K4").Value = 7000
$ws.Range("L4").Value = 7500
$ws.Range("M4").Value = 7200
$ws.Range("P4").Value = 720
$ws.Range("D5").Value = 45037
$ws.Range("J5").Value = 40
$ws.Range("K5").Value = 9000
$ws.Range("L5").Value = 10000
$ws.Range("M5").Value = 9500
$ws.Range("P5").Value = 950
$ws.Range("D6").Value = 44308
$ws.Range("J6").Value = 100
$ws.Range("K6").Value = 5000
$ws.Range("L6").Value = 5500
$ws.Range("M6").Value = 5250
$ws.Range("P6").Value = 525
$ws.Range("D7").Value = 44797
$ws.Range("J7").Value = 50
$ws.Range("K7").Value = 6500
$ws.Range("L7").Value = 7000
$ws.Range("M7").Value = 6700
$ws.Range("P7").Value = 670
$ws.Range("D8").Value = 44565
$ws.Range("J8").Value = 40
$ws.Range("K8").Value = 12000
$ws.Range("L8").Value = 13000
$ws.Range("M8").Value = 12500
$ws.Range("P8").Value = 1250
$ws.Range("D9").Value = 44727
$ws.Range("J9").Value = 130
$ws.Range("K9").Value = 7000
$ws.Range("L9").Value = 7500
$ws.Range("M9").Value = 7192
$ws.Range("P9").Value = 719
$ws.Range("D10").Value = 44194
$ws.Range("J10").Value = 100
$ws.Range("L10").Value = 9000
$ws.Range("M10").Value = 8500
$ws.Range("P10").Value = 850
$ws.Range("D11").Value = 44230
$ws.Range("J11").Value = 60
$ws.Range("K11").Value = 9000
$ws.Range("L11").Value = 10000
$ws.Range("M11").Value = 9500
$ws.Range("P11").Value = 950
$ws.Range("D12").Value = 44680
$ws.Range("J12").Value = 50
$ws.Range("K12").Value = 6000
$ws.Range("L12").Value = 6500
$ws.Range("M12").Value = 6300
$ws.Range("P12").Value = 630
$ws.Range("D13").Value = 44967
$ws.Range("J13").Value = 40
$ws.Range("K13").Value = 7000
$ws.Range("L13").Value = 7500
$ws.Range("M13").Value = 7250
$ws.Range("P13").Value = 725
$ws.Range("D14").Value = 44651
$ws.Range("J14").Value = 250
$ws.Range("K14").Value = 6000
$ws.Range("L14").Value = 6500
$ws.Range("M14").Value = 6200
$ws.Range("P14").Value = 620
$ws.Range("D15").Value = 44763
$ws.Range("J15").Value = 50
$ws.Range("M15").Value = 9600
$ws.Range("P15").Value = 960
$ws.Range("D16").Value = 44923
$ws.Range("K16").Value = 7000
$ws.Range("L16").Value = 7500
$ws.Range("M16").Value = 7250
$ws.Range("P16").Value = 725
$ws.Range("D17").Value = 44923
$ws.Range("K17").Value = 7000
$ws.Range("L17").Value = 7500
$ws.Range("M17").Value = 7250
$ws.Range("P17").Value = 725
$ws.Range("D18").Value = 44453
$ws.Range("J18").Value = 100
$ws.Range("K18").Value = 6500
$ws.Range("L18").Value = 7000
$ws.Range("M18").Value = 6750
$ws.Range("P18").Value = 675
$ws.Range("D19").Value = 44755
$ws.Range("J19").Value = 40
$ws.Range("K19").Value = 9000
$ws.Range("L19").Value = 10000
$ws.Range("M19").Value = 9500
$ws.Range("P19").Value = 950
$ws.Range("D20").Value = 44582
$ws.Range("K20").Value = 8000
$ws.Range("L20").Value = 8500
$ws.Range("M20").Value = 8200
$ws.Range("P20").Value = 820
$ws.Range("D21").Value = 44383
$ws.Range("J21").Value = 60
$ws.Range("K21").Value = 7500
$ws.Range("L21").Value = 8000
$ws.Range("M21").Value = 7750
$ws.Range("P21").Value = 775
$ws.Range("D22").Value = 44719
$ws.Range("K22").Value = 7000
$ws.Range("L22").Value = 7500
$ws.Range("M22").Value = 7250
$ws.Range("P22").Value = 725
$ws.Range("D23").Value = 44813
$ws.Range("J23").Value = 50
$ws.Range("K23").Value = 6500
$ws.Range("L23").Value = 7000
$ws.Range("M23").Value = 6700
$ws.Range("P23").Value = 670
$ws.Range("D24").Value = 45020
$ws.Range("J24").Value = 40
$ws.Range("K24").Value = 9000
$ws.Range("L24").Value = 10000
$ws.Range("M24").Value = 9500
$ws.Range("P24").Value = 950
$ws.Range("D25").Value = 44995
$ws.Range("J25").Value = 40
$ws.Range("K25").Value = 7500
$ws.Range("L25").Value = 8000
$ws.Range("M25").Value = 7750
$ws.Range("P25").Value = 775
$ws.Range("D26").Value = 44351
$ws.Range("J26").Value = 50
$ws.Range("K26").Value = 6000
$ws.Range("L26").Value = 6500
$ws.Range("M26").Value = 6300
$ws.Range("P26").Value = 630
$ws.Range("D27").Value = 44838
$ws.Range("K27").Value = 6500
$ws.Range("L27").Value = 7000
$ws.Range("M27").Value = 6700
$ws.Range("P27").Value = 670
$ws.Range("D28").Value = 44496
$ws.Range("J28").Value = 350
$ws.Range("K28").Value = 650
$ws.Range("L28").Value = 700
$ws.Range("M28").Value = 671
$ws.Range("N28").Value = '$/atado 0,5 a 1 kilo'
$ws.Range("O28").Value = 'Región de Ñuble'
$ws.Range("P28").Value = 671
$ws.Range("Q28").Value = 1
$ws.Range("D29").Value = 44421
$ws.Range("K29").Value = 7000
$ws.Range("L29").Value = 7500
$ws.Range("M29").Value = 7250
$ws.Range("P29").Value = 725
$ws.Range("D30").Value = 44357
$ws.Range("J30").Value = 50
$ws.Range("K30").Value = 6000
$ws.Range("L30").Value = 6500
$ws.Range("M30").Value = 6200
$ws.Range("P30").Value = 620
$ws.Range("D31").Value = 44349
$ws.Range("J31").Value = 60
$ws.Range("M31").Value = 6250
$ws.Range("P31").Value = 625
$ws.Range("D32").Value = 44635
$ws.Range("J32").Value = 170
$ws.Range("K32").Value = 7500
$ws.Range("L32").Value = 8000
$ws.Range("M32").Value = 7765
$ws.Range("P32").Value = 776
$ws.Range("D33").Value = 44523
$ws.Range("J33").Value = 50
$ws.Range("K33").Value = 6500
$ws.Range("L33").Value = 7000
$ws.Range("M33").Value = 6800
$ws.Range("P33").Value = 680
$ws.Range("D34").Value = 44810
$ws.Range("J34").Value = 40
$ws.Range("K34").Value = 6500
$ws.Range("L34").Value = 7000
$ws.Range("M34").Value = 6750
$ws.Range("P34").Value = 675
$ws.Range("D35").Value = 44365
$ws.Range("J35").Value = 50
$ws.Range("K35").Value = 6000
$ws.Range("L35").Value = 6500
$ws.Range("M35").Value = 6200
$ws.Range("P35").Value = 620
$ws.Range("D36").Value = 44355
$ws.Range("J36").Value = 50
$ws.Range("K36").Value = 6000
$ws.Range("L36").Value = 6500
$ws.Range("M36").Value = 6300
$ws.Range("P36").Value = 630
$ws.Range("D37").Value = 45044
$ws.Range("J37").Value = 220
$ws.Range("K37").Value = 7000
$ws.Range("L37").Value = 8000
$ws.Range("M37").Value = 7545
$ws.Range("P37").Value = 754
$ws.Range("D38").Value = 44657
$ws.Range("J38").Value = 220
$ws.Range("K38").Value = 6000
$ws.Range("L38").Value = 6500
$ws.Range("M38").Value = 6273
$ws.Range("P38").Value = 627
$ws.Range("D39").Value = 44881
$ws.Range("J39").Value = 50
$ws.Range("K39").Value = 13000
$ws.Range("L39").Value = 14000
$ws.Range("M39").Value = 13400
$ws.Range("P39").Value = 1340
$ws.Range("D40").Value = 44784
$ws.Range("J40").Value = 40
$ws.Range("D41").Value = 45022
$ws.Range("J41").Value = 40
$ws.Range("K41").Value = 9000
$ws.Range("L41").Value = 10000
$ws.Range("M41").Value = 9500
$ws.Range("P41").Value = 950
$ws.Range("D42").Value = 44428
$ws.Range("J42").Value = 50
$ws.Range("M42").Value = 7800
$ws.Range("P42").Value = 780
$ws.Range("D43").Value = 44334
$ws.Range("J43").Value = 60
$ws.Range("K43").Value = 6500
$ws.Range("L43").Value = 7000
$ws.Range("M43").Value = 6750
$ws.Range("P43").Value = 675
$ws.Range("D44").Value = 45049
$ws.Range("J44").Value = 90
$ws.Range("K44").Value = 8500
$ws.Range("L44").Value = 10000
$ws.Range("M44").Value = 9222
$ws.Range("N44").Value = '$/cuna 10 kilos'
$ws.Range("O44").Value = 'Región Metropolitana'
$ws.Range("P44").Value = 922
$ws.Range("Q44").Value = 10
$ws.Range("D46").Value = 44993
$ws.Range("J46").Value = 40
$ws.Range("K46").Value = 7500
$ws.Range("L46").Value = 8000
$ws.Range("M46").Value = 7750
$ws.Range("P46").Value = 775
$ws.Range("D47").Value = 44362
$ws.Range("K47").Value = 6000
$ws.Range("L47").Value = 6500
$ws.Range("M47").Value = 6300
$ws.Range("P47").Value = 630
$ws.Range("D48").Value = 44630
$ws.Range("J48").Value = 150
$ws.Range("K48").Value = 10000
$ws.Range("L48").Value = 10000
$ws.Range("M48").Value = 10000
$ws.Range("P48").Value = 1000
$ws.Range("D49").Value = 44670
$ws.Range("J49").Value = 170
$ws.Range("K49").Value = 5000
$ws.Range("L49").Value = 6000
$ws.Range("M49").Value = 5471
$ws.Range("P49").Value = 547
$ws.Range("D50").Value = 44299
$ws.Range("J50").Value = 100
$ws.Range("K50").Value = 8000
$ws.Range("L50").Value = 9000
$ws.Range("M50").Value = 8500
$ws.Range("P50").Value = 850
$ws.Range("D51").Value = 44708
$ws.Range("J51").Value = 100
$ws.Range("K51").Value = 5500
$ws.Range("L51").Value = 6000
$ws.Range("M51").Value = 5750
$ws.Range("P51").Value = 575
$ws.Range("D52").Value = 45062
$ws.Range("J52").Value = 50
$ws.Range("K52").Value = 10000
$ws.Range("L52").Value = 11000
$ws.Range("M52").Value = 10400
$ws.Range("P52").Value = 1040
$ws.Range("D53").Value = 44757
$ws.Range("J53").Value = 40
$ws.Range("K53").Value = 8000
$ws.Range("L53").Value = 8500
$ws.Range("M53").Value = 8250
$ws.Range("P53").Value = 825
$ws.Range("D54").Value = 44327
$ws.Range("J54").Value = 60
$ws.Range("M54").Value = 7250
$ws.Range("P54").Value = 725
$ws.Range("D55").Value = 44476
$ws.Range("J55").Value = 80
$ws.Range("K55").Value = 5000
$ws.Range("L55").Value = 5500
$ws.Range("M55").Value = 5312
$ws.Range("P55").Value = 531
$ws.Range("D56").Value = 44433
$ws.Range("J56").Value = 100
$ws.Range("K56").Value = 7000
$ws.Range("L56").Value = 7500
$ws.Range("M56").Value = 7250
$ws.Range("P56").Value = 725
$ws.Range("D57").Value = 44328
$ws.Range("J57").Value = 60
$ws.Range("K57").Value = 7000
$ws.Range("L57").Value = 7500
$ws.Range("M57").Value = 7250
$ws.Range("P57").Value = 725
$ws.Range("D58").Value = 45071
$ws.Range("J58").Value = 220
$ws.Range("K58").Value = 8000
$ws.Range("L58").Value = 9000
$ws.Range("M58").Value = 8545
$ws.Range("P58").Value = 854
$ws.Range("D59").Value = 44188
$ws.Range("J59").Value = 80
$ws.Range("K59").Value = 8000
$ws.Range("L59").Value = 8500
$ws.Range("M59").Value = 8250
$ws.Range("P59").Value = 825
$ws.Range("D60").Value = 45035
$ws.Range("J60").Value = 40
$ws.Range("K60").Value = 10000
$ws.Range("L60").Value = 11000
$ws.Range("M60").Value = 10500
$ws.Range("P60").Value = 1050
$ws.Range("D61").Value = 45135
$ws.Range("K61").Value = 6500
$ws.Range("L61").Value = 7000
$ws.Range("M61").Value = 6750
$ws.Range("P61").Value = 675
$ws.Range("D62").Value = 45030
$ws.Range("J62").Value = 110
$ws.Range("K62").Value = 11000
$ws.Range("L62").Value = 12000
$ws.Range("M62").Value = 11455
$ws.Range("P62").Value = 1146
$ws.Range("D63").Value = 44316
$ws.Range("J63").Value = 100
$ws.Range("K63").Value = 6000
$ws.Range("L63").Value = 6500
$ws.Range("M63").Value = 6250
$ws.Range("P63").Value = 625
$ws.Range("D64").Value = 44225
$ws.Range("J64").Value = 60
$ws.Range("K64").Value = 7500
$ws.Range("L64").Value = 8000
$ws.Range("M64").Value = 7750
$ws.Range("P64").Value = 775
$ws.Range("D65").Value = 44376
$ws.Range("J65").Value = 100
$ws.Range("K65").Value = 6000
$ws.Range("L65").Value = 6500
$ws.Range("M65").Value = 6250
$ws.Range("P65").Value = 625
$ws.Range("D66").Value = 44687
$ws.Range("J66").Value = 150
$ws.Range("K66").Value = 6500
$ws.Range("L66").Value = 7000
$ws.Range("M66").Value = 6733
$ws.Range("P66").Value = 673
$ws.Range("D67").Value = 44678
$ws.Range("J67").Value = 40
$ws.Range("K67").Value = 6000
$ws.Range("L67").Value = 6500
$ws.Range("M67").Value = 6250
$ws.Range("P67").Value = 625
$ws.Range("D68").Value = 45042
$ws.Range("J68").Value = 140
$ws.Range("K68").Value = 7000
$ws.Range("L68").Value = 7500
$ws.Range("M68").Value = 7214
$ws.Range("P68").Value = 721
$ws.Range("D69").Value = 44966
$ws.Range("J69").Value = 40
$ws.Range("L69").Value = 7500
$ws.Range("M69").Value = 7250
$ws.Range("P69").Value = 725
$ws.Range("D70").Value = 44204
$ws.Range("J70").Value = 80
$ws.Range("K70").Value = 7000
$ws.Range("L70").Value = 7500
$ws.Range("M70").Value = 7188
$ws.Range("P70").Value = 719
$ws.Range("D71").Value = 44855
$ws.Range("J71").Value = 50
$ws.Range("K71").Value = 8000
$ws.Range("L71").Value = 8500
$ws.Range("M71").Value = 8200
$ws.Range("P71").Value = 820
$ws.Range("D72").Value = 45079
$ws.Range("J72").Value = 40
$ws.Range("K72").Value = 5000
$ws.Range("L72").Value = 5500
$ws.Range("M72").Value = 5250
$ws.Range("P72").Value = 525
$ws.Range("D73").Value = 45099
$ws.Range("J73").Value = 220
$ws.Range("K73").Value = 6500
$ws.Range("L73").Value = 7000
$ws.Range("M73").Value = 6727
$ws.Range("P73").Value = 673
$ws.Range("D74").Value = 44405
$ws.Range("J74").Value = 80
$ws.Range("K74").Value = 7500
$ws.Range("L74").Value = 8000
$ws.Range("M74").Value = 7688
$ws.Range("P74").Value = 769
$ws.Range("D75").Value = 44273
$ws.Range("J75").Value = 80
$ws.Range("K75").Value = 7000
$ws.Range("L75").Value = 8000
$ws.Range("M75").Value = 7500
$ws.Range("P75").Value = 750
$ws.Range("D76").Value = 44897
$ws.Range("J76").Value = 100
$ws.Range("K76").Value = 13000
$ws.Range("L76").Value = 14000
$ws.Range("M76").Value = 13500
$ws.Range("P76").Value = 1350
$ws.Range("D77").Value = 45093
$ws.Range("J77").Value = 40
$ws.Range("K77").Value = 6500
$ws.Range("L77").Value = 7000
$ws.Range("M77").Value = 6750
$ws.Range("P77").Value = 675
$ws.Range("D78").Value = 44761
$ws.Range("J78").Value = 45
$ws.Range("K78").Value = 10000
$ws.Range("L78").Value = 11000
$ws.Range("M78").Value = 10333
$ws.Range("P78").Value = 1033
$ws.Range("D79").Value = 44336
$ws.Range("J79").Value = 60
$ws.Range("M79").Value = 6250
$ws.Range("P79").Value = 625
$ws.Range("D80").Value = 44908
$ws.Range("J80").Value = 40
$ws.Range("K80").Value = 9000
$ws.Range("L80").Value = 10000
$ws.Range("M80").Value = 9500
$ws.Range("P80").Value = 950
$ws.Range("D81").Value = 44238
$ws.Range("J81").Value = 100
$ws.Range("K81").Value = 8000
$ws.Range("L81").Value = 8500
$ws.Range("M81").Value = 8250
$ws.Range("P81").Value = 825
$ws.Range("D82").Value = 44474
$ws.Range("J82").Value = 50
$ws.Range("K82").Value = 6000
$ws.Range("L82").Value = 6500
$ws.Range("M82").Value = 6300
$ws.Range("P82").Value = 630
$ws.Range("D83").Value = 44782
$ws.Range("J83").Value = 100
$ws.Range("K83").Value = 7000
$ws.Range("L83").Value = 7500
$ws.Range("M83").Value = 7250
$ws.Range("P83").Value = 725
$ws.Range("D84").Value = 44776
$ws.Range("J84").Value = 150
$ws.Range("K84").Value = 6500
$ws.Range("L84").Value = 7000
$ws.Range("M84").Value = 6733
$ws.Range("P84").Value = 673
$ws.Range("D85").Value = 44771
$ws.Range("J85").Value = 40
$ws.Range("K85").Value = 7000
$ws.Range("L85").Value = 7500
$ws.Range("M85").Value = 7250
$ws.Range("P85").Value = 725
$ws.Range("D86").Value = 44706
$ws.Range("J86").Value = 100
$ws.Range("K86").Value = 5500
$ws.Range("L86").Value = 6000
$ws.Range("M86").Value = 5750
$ws.Range("P86").Value = 575
$ws.Range("D87").Value = 44160
$ws.Range("J87").Value = 100
$ws.Range("K87").Value = 9000
$ws.Range("L87").Value = 9500
$ws.Range("M87").Value = 9250
$ws.Range("P87").Value = 925
$ws.Range("D88").Value = 45028
$ws.Range("J88").Value = 160
$ws.Range("L88").Value = 9000
$ws.Range("M88").Value = 8375
$ws.Range("P88").Value = 838
$ws.Range("D89").Value = 44399
$ws.Range("J89").Value = 60
$ws.Range("K89").Value = 9000
$ws.Range("L89").Value = 10000
$ws.Range("M89").Value = 9500
$ws.Range("P89").Value = 950
$ws.Range("D90").Value = 44981
$ws.Range("J90").Value = 50
$ws.Range("K90").Value = 7000
$ws.Range("L90").Value = 7500
$ws.Range("M90").Value = 7200
$ws.Range("P90").Value = 720
$ws.Range("D91").Value = 44659
$ws.Range("J91").Value = 250
$ws.Range("K91").Value = 6000
$ws.Range("L91").Value = 6500
$ws.Range("M91").Value = 6200
$ws.Range("P91").Value = 620
$ws.Range("D92").Value = 44509
$ws.Range("J92").Value = 80
$ws.Range("K92").Value = 6500
$ws.Range("L92").Value = 7000
$ws.Range("M92").Value = 6750
$ws.Range("P92").Value = 675
$ws.Range("D93").Value = 44811
$ws.Range("J93").Value = 50
$ws.Range("K93").Value = 5500
$ws.Range("L93").Value = 6000
$ws.Range("M93").Value = 5800
$ws.Range("P93").Value = 580
$ws.Range("D94").Value = 44699
$ws.Range("J94").Value = 50
$ws.Range("K94").Value = 5500
$ws.Range("L94").Value = 6000
$ws.Range("M94").Value = 5700
$ws.Range("P94").Value = 570
$ws.Range("D95").Value = 44951
$ws.Range("K95").Value = 7000
$ws.Range("L95").Value = 7500
$ws.Range("M95").Value = 7300
$ws.Range("P95").Value = 730
$ws.Range("D96").Value = 44552
$ws.Range("J96").Value = 60
$ws.Range("K96").Value = 11000
$ws.Range("L96").Value = 12000
$ws.Range("M96").Value = 11500
$ws.Range("P96").Value = 1150
$ws.Range("D97").Value = 44575
$ws.Range("J97").Value = 50
$ws.Range("K97").Value = 11000
$ws.Range("L97").Value = 12000
$ws.Range("M97").Value = 11600
$ws.Range("P97").Value = 1160
$ws.Range("D98").Value = 44747
$ws.Range("J98").Value = 150
$ws.Range("K98").Value = 8000
$ws.Range("L98").Value = 8500
$ws.Range("M98").Value = 8233
$ws.Range("P98").Value = 823
$ws.Range("D99").Value = 44636
$ws.Range("J99").Value = 220
$ws.Range("K99").Value = 8000
$ws.Range("L99").Value = 9000
$ws.Range("M99").Value = 8545
$ws.Range("P99").Value = 854
$ws.Range("D100").Value = 44580
$ws.Range("K100").Value = 10000
$ws.Range("L100").Value = 11000
$ws.Range("M100").Value = 10500
$ws.Range("P100").Value = 1050
$ws.Range("D101").Value = 45091
$ws.Range("J101").Value = 50
$ws.Range("K101").Value = 6500
$ws.Range("L101").Value = 7000
$ws.Range("M101").Value = 6700
$ws.Range("P101").Value = 670
$ws.Range("D102").Value = 44512
$ws.Range("J102").Value = 60
$ws.Range("D104").Value = 44246
$ws.Range("J104").Value = 60
$ws.Range("K104").Value = 9000
$ws.Range("L104").Value = 10000
$ws.Range("M104").Value = 9500
$ws.Range("P104").Value = 950
$ws.Range("D105").Value = 44665
$ws.Range("J105").Value = 100
$ws.Range("K105").Value = 6500
$ws.Range("L105").Value = 7000
$ws.Range("M105").Value = 6750
$ws.Range("P105").Value = 675
$ws.Range("D106").Value = 45007
$ws.Range("J106").Value = 50
$ws.Range("K106").Value = 6000
$ws.Range("L106").Value = 6500
$ws.Range("M106").Value = 6200
$ws.Range("P106").Value = 620
$ws.Range("D107").Value = 44209
$ws.Range("J107").Value = 80
$ws.Range("K107").Value = 7500
$ws.Range("L107").Value = 8000
$ws.Range("M107").Value = 7688
$ws.Range("P107").Value = 769
$ws.Range("D108").Value = 44435
$ws.Range("J108").Value = 100
$ws.Range("M108").Value = 7250
$ws.Range("P108").Value = 725
$ws.Range("D109").Value = 44425
$ws.Range("J109").Value = 60
$ws.Range("K109").Value = 6500
$ws.Range("L109").Value = 7000
$ws.Range("M109").Value = 6750
$ws.Range("P109").Value = 675
$ws.Range("D111").Value = 44741
$ws.Range("J111").Value = 100
$ws.Range("K111").Value = 8000
$ws.Range("L111").Value = 8500
$ws.Range("M111").Value = 8250
$ws.Range("P111").Value = 825
$ws.Range("D112").Value = 44972
$ws.Range("J112").Value = 50
$ws.Range("K112").Value = 7000
$ws.Range("L112").Value = 7500
$ws.Range("M112").Value = 7200
$ws.Range("P112").Value = 720
$ws.Range("D113").Value = 44342
$ws.Range("K113").Value = 6000
$ws.Range("L113").Value = 6500
$ws.Range("M113").Value = 6300
$ws.Range("P113").Value = 630
$ws.Range("D114").Value = 44489
$ws.Range("J114").Value = 50
$ws.Range("K114").Value = 6000
$ws.Range("L114").Value = 6500
$ws.Range("M114").Value = 6300
$ws.Range("P114").Value = 630
$ws.Range("D115").Value = 44469
$ws.Range("K115").Value = 6000
$ws.Range("L115").Value = 6500
$ws.Range("M115").Value = 6250
$ws.Range("P115").Value = 625
$ws.Range("D116").Value = 44607
$ws.Range("J116").Value = 100
$ws.Range("K116").Value = 12000
$ws.Range("L116").Value = 13000
$ws.Range("M116").Value = 12500
$ws.Range("P116").Value = 1250
$ws.Range("D117").Value = 45106
$ws.Range("J117").Value = 50
$ws.Range("K117").Value = 7000
$ws.Range("L117").Value = 7500
$ws.Range("M117").Value = 7200
$ws.Range("P117").Value = 720
$ws.Range("D118").Value = 45015
$ws.Range("J118").Value = 150
$ws.Range("K118").Value = 7500
$ws.Range("L118").Value = 8000
$ws.Range("M118").Value = 7667
$ws.Range("P118").Value = 767
$ws.Range("D119").Value = 44692
$ws.Range("D120").Value = 44166
$ws.Range("J120").Value = 100
$ws.Range("K120").Value = 8000
$ws.Range("L120").Value = 9000
$ws.Range("M120").Value = 8500
$ws.Range("P120").Value = 850
$ws.Range("D121").Value = 45077
$ws.Range("J121").Value = 40
$ws.Range("K121").Value = 5000
$ws.Range("L121").Value = 5500
$ws.Range("M121").Value = 5250
$ws.Range("O121").Value = 'Región Metropolitana'
$ws.Range("P121").Value = 525
$ws.Range("D122").Value = 44231
$ws.Range("J122").Value = 70
$ws.Range("K122").Value = 7500
$ws.Range("L122").Value = 8000
$ws.Range("M122").Value = 7714
$ws.Range("P122").Value = 771
$ws.Range("D123").Value = 44292
$ws.Range("J123").Value = 50
$ws.Range("K123").Value = 10000
$ws.Range("L123").Value = 11000
$ws.Range("M123").Value = 10600
$ws.Range("P123").Value = 1060
$ws.Range("D124").Value = 45119
$ws.Range("J124").Value = 70
$ws.Range("K124").Value = 7500
$ws.Range("L124").Value = 8000
$ws.Range("M124").Value = 7643
$ws.Range("P124").Value = 764
$ws.Range("D125").Value = 44526
$ws.Range("J125").Value = 40
$ws.Range("K125").Value = 8000
$ws.Range("L125").Value = 8500
$ws.Range("M125").Value = 8250
$ws.Range("P125").Value = 825
$ws.Range("D126").Value = 44631
$ws.Range("K126").Value = 9000
$ws.Range("L126").Value = 9500
$ws.Range("M126").Value = 9273
$ws.Range("O126").Value = 'Provincia de Chacabuco'
$ws.Range("P126").Value = 927
$ws.Range("D127").Value = 44645
$ws.Range("J127").Value = 160
$ws.Range("K127").Value = 6000
$ws.Range("L127").Value = 6500
$ws.Range("M127").Value = 6250
$ws.Range("P127").Value = 625
$ws.Range("D128").Value = 44264
$ws.Range("K128").Value = 8000
$ws.Range("L128").Value = 8500
$ws.Range("M128").Value = 8200
$ws.Range("P128").Value = 820
$ws.Range("D129").Value = 44644
$ws.Range("J129").Value = 220
$ws.Range("K129").Value = 6500
$ws.Range("L129").Value = 7000
$ws.Range("M129").Value = 6727
$ws.Range("P129").Value = 673
$ws.Range("D130").Value = 44313
$ws.Range("J130").Value = 60
$ws.Range("K130").Value = 6000
$ws.Range("L130").Value = 6500
$ws.Range("M130").Value = 6250
$ws.Range("P130").Value = 625
$ws.Range("D131").Value = 44979
$ws.Range("J131").Value = 40
$ws.Range("K131").Value = 7500
$ws.Range("L131").Value = 8000
$ws.Range("M131").Value = 7750
$ws.Range("P131").Value = 775
$ws.Range("D132").Value = 45133
$ws.Range("J132").Value = 50
$ws.Range("K132").Value = 7000
$ws.Range("L132").Value = 7500
$ws.Range("M132").Value = 7200
$ws.Range("P132").Value = 720
